# Update market/profit data cells (H-N columns) across several leve-profit
# sheets, as produced by the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1224.6666
$ws.Range("I12").Value = 975
$ws.Range("J12").Value = 1349.5
$ws.Range("K12").Value = 975
$ws.Range("L12").Value = 1349.5
$ws.Range("M12").Value = -805
$ws.Range("N12").Value = -1689.5

$ws.Range("H29").Value = 519.8
$ws.Range("I29").Value = 519.8
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 1559.4
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -1278.4
$ws.Range("N29").ClearContents()

$ws.Range("H38").Value = 153.5
$ws.Range("I38").Value = 153.5
$ws.Range("K38").Value = 460.5
$ws.Range("M38").Value = -88.5

$ws.Range("H58").Value = 46881052
$ws.Range("J58").Value = 100019000
$ws.Range("L58").Value = 300057000
$ws.Range("N58").Value = -300057300

$ws.Range("H64").Value = 40007096
$ws.Range("I64").Value = 58830550
$ws.Range("J64").Value = 7249.875
$ws.Range("K64").Value = 58830550
$ws.Range("L64").Value = 7249.875
$ws.Range("M64").Value = -58830302
$ws.Range("N64").Value = -7745.875

$ws.Range("H67").Value = 40007096
$ws.Range("I67").Value = 58830550
$ws.Range("J67").Value = 7249.875
$ws.Range("K67").Value = 58830550
$ws.Range("L67").Value = 7249.875
$ws.Range("M67").Value = -58829692
$ws.Range("N67").Value = -8965.875

$ws.Range("H87").Value = 72000
$ws.Range("J87").Value = 72000
$ws.Range("L87").Value = 72000
$ws.Range("N87").Value = -74496

$ws.Range("H90").Value = 72000
$ws.Range("J90").Value = 72000
$ws.Range("L90").Value = 216000
$ws.Range("N90").Value = -228480

$ws.Range("H132").Value = 1835.48
$ws.Range("I132").Value = 1944.409
$ws.Range("J132").Value = 1036.6666
$ws.Range("K132").Value = 5833.227000000001
$ws.Range("L132").Value = 3109.9998
$ws.Range("M132").Value = -3303.227000000001
$ws.Range("N132").Value = -8169.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2416.6667
$ws.Range("I45").Value = 1321.5
$ws.Range("K45").Value = 1321.5
$ws.Range("M45").Value = -944.5

$ws.Range("H61").Value = 8812.412
$ws.Range("I61").Value = 1788.1666
$ws.Range("J61").Value = 12643.818
$ws.Range("K61").Value = 1788.1666
$ws.Range("L61").Value = 12643.818
$ws.Range("M61").Value = -1576.1666
$ws.Range("N61").Value = -13067.818

$ws.Range("H74").Value = 66650.64
$ws.Range("I74").Value = 95308.94
$ws.Range("K74").Value = 95308.94
$ws.Range("M74").Value = -94434.94

$ws.Range("H77").Value = 66650.64
$ws.Range("I77").Value = 95308.94
$ws.Range("K77").Value = 476544.7
$ws.Range("M77").Value = -472176.7

$ws.Range("H97").Value = 4631729.5
$ws.Range("I97").Value = 2575.75
$ws.Range("K97").Value = 2575.75
$ws.Range("M97").Value = -2079.75

$ws.Range("H132").Value = 5112.036
$ws.Range("I132").Value = 2384.8286
$ws.Range("K132").Value = 7154.485799999999
$ws.Range("M132").Value = -4624.485799999999

$ws.Range("H135").Value = 80390
$ws.Range("J135").Value = 80390
$ws.Range("L135").Value = 80390
$ws.Range("N135").Value = -90530

$ws.Range("H136").Value = 8812.412
$ws.Range("I136").Value = 1788.1666
$ws.Range("J136").Value = 12643.818
$ws.Range("K136").Value = 5364.4998
$ws.Range("L136").Value = 37931.454
$ws.Range("M136").Value = -2814.4998
$ws.Range("N136").Value = -43031.454

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2331.0435
$ws.Range("I94").Value = 1845.6923
$ws.Range("J94").Value = 2962
$ws.Range("K94").Value = 1845.6923
$ws.Range("L94").Value = 2962
$ws.Range("M94").Value = -1394.6923
$ws.Range("N94").Value = -3864

$ws.Range("H134").Value = 8935966
$ws.Range("I134").Value = 22729402
$ws.Range("J134").Value = 10802.529
$ws.Range("K134").Value = 68188206
$ws.Range("L134").Value = 32407.587
$ws.Range("M134").Value = -68185671
$ws.Range("N134").Value = -37477.587

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8076.543
$ws.Range("I31").Value = 2894.4285
$ws.Range("J31").Value = 11531.286
$ws.Range("K31").Value = 2894.4285
$ws.Range("L31").Value = 11531.286
$ws.Range("M31").Value = -2599.4285
$ws.Range("N31").Value = -12121.286

$ws.Range("H34").Value = 8076.543
$ws.Range("I34").Value = 2894.4285
$ws.Range("J34").Value = 11531.286
$ws.Range("K34").Value = 2894.4285
$ws.Range("L34").Value = 11531.286
$ws.Range("M34").Value = -2692.4285
$ws.Range("N34").Value = -11935.286

$ws.Range("H132").Value = 5227.222
$ws.Range("I132").Value = 3430
$ws.Range("K132").Value = 10290
$ws.Range("M132").Value = -7760

$ws.Range("H134").Value = 4967.372
$ws.Range("I134").Value = 1922.875
$ws.Range("K134").Value = 5768.625
$ws.Range("M134").Value = -3233.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 12292381
$ws.Range("I4").Value = 23500090
$ws.Range("J4").Value = 1084672.6
$ws.Range("K4").Value = 70500270
$ws.Range("L4").Value = 3254017.8
$ws.Range("M4").Value = -70500158
$ws.Range("N4").Value = -3254241.8

$ws.Range("H34").Value = 6124.9316
$ws.Range("I34").Value = 9000
$ws.Range("J34").Value = 6058.07
$ws.Range("K34").Value = 27000
$ws.Range("L34").Value = 18174.21
$ws.Range("M34").Value = -26916
$ws.Range("N34").Value = -18342.21

$ws.Range("H55").Value = 5890758.5
$ws.Range("J55").Value = 6258807
$ws.Range("L55").Value = 18776421
$ws.Range("N55").Value = -18776775

$ws.Range("H56").Value = 6999
$ws.Range("I56").Value = 6999
$ws.Range("K56").Value = 6999
$ws.Range("M56").Value = -6469

$ws.Range("H98").Value = 2349.7693
$ws.Range("J98").Value = 4282
$ws.Range("L98").Value = 12846
$ws.Range("N98").Value = -15842

$ws.Range("H117").Value = 973.4
$ws.Range("I117").Value = 919.5714
$ws.Range("J117").Value = 1020.5
$ws.Range("K117").Value = 2758.7142
$ws.Range("L117").Value = 3061.5
$ws.Range("M117").Value = 683.2857999999997
$ws.Range("N117").Value = -9945.5

$ws.Range("H121").Value = 1364.125
$ws.Range("I121").Value = 1285.9166
$ws.Range("J121").Value = 1598.75
$ws.Range("K121").Value = 3857.7498
$ws.Range("L121").Value = 4796.25
$ws.Range("M121").Value = -2547.7498
$ws.Range("N121").Value = -7416.25

$ws.Range("H122").Value = 3143594
$ws.Range("I122").Value = 3536493
$ws.Range("J122").Value = 404
$ws.Range("K122").Value = 31828437
$ws.Range("L122").Value = 3636
$ws.Range("M122").Value = -31825987
$ws.Range("N122").Value = -8536

$ws.Range("H134").Value = 58303.684
$ws.Range("I134").Value = 104049.5
$ws.Range("K134").Value = 312148.5
$ws.Range("M134").Value = -307078.5

$ws.Range("H139").Value = 56981.895
$ws.Range("I139").Value = 65788.94
$ws.Range("K139").Value = 197366.82
$ws.Range("M139").Value = -192226.82

$ws.Range("H140").Value = 71048.31
$ws.Range("I140").Value = 88104.61
$ws.Range("J140").Value = 5665.8335
$ws.Range("K140").Value = 264313.83
$ws.Range("L140").Value = 16997.5005
$ws.Range("M140").Value = -259133.83
$ws.Range("N140").Value = -27357.5005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1960222.2
$ws.Range("I122").Value = 2589048.8
$ws.Range("J122").Value = 3873.5557
$ws.Range("K122").Value = 7767146.399999999
$ws.Range("L122").Value = 11620.6671
$ws.Range("M122").Value = -7764696.399999999
$ws.Range("N122").Value = -16520.6671

$ws.Range("H132").Value = 2296
$ws.Range("I132").Value = 2225.423
$ws.Range("K132").Value = 6676.268999999999
$ws.Range("M132").Value = -4146.268999999999

$ws.Range("H133").Value = 80000
$ws.Range("J133").Value = 80000
$ws.Range("L133").Value = 80000
$ws.Range("N133").Value = -90120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4344.8066
$ws.Range("I7").Value = 2399.389
$ws.Range("K7").Value = 2399.389
$ws.Range("M7").Value = -2287.389

$ws.Range("H22").Value = 5453.174
$ws.Range("I22").Value = 594.2
$ws.Range("J22").Value = 14563.75
$ws.Range("K22").Value = 594.2
$ws.Range("L22").Value = 14563.75
$ws.Range("M22").Value = -299.2
$ws.Range("N22").Value = -15153.75

$ws.Range("H27").Value = 5453.174
$ws.Range("I27").Value = 594.2
$ws.Range("J27").Value = 14563.75
$ws.Range("K27").Value = 594.2
$ws.Range("L27").Value = 14563.75
$ws.Range("M27").Value = -487.2
$ws.Range("N27").Value = -14777.75

$ws.Range("H55").Value = 891.75
$ws.Range("I55").Value = 670
$ws.Range("K55").Value = 670
$ws.Range("M55").Value = -497

$ws.Range("H101").Value = 53247
$ws.Range("J101").Value = 53247
$ws.Range("L101").Value = 53247
$ws.Range("N101").Value = -59737

$ws.Range("H104").Value = 35023.668
$ws.Range("J104").Value = 35023.668
$ws.Range("L104").Value = 35023.668
$ws.Range("N104").Value = -42011.668

$ws.Range("H126").Value = 4344.8066
$ws.Range("I126").Value = 2399.389
$ws.Range("K126").Value = 7198.167
$ws.Range("M126").Value = -4728.167

$ws.Range("H136").Value = 8055.222
$ws.Range("I136").Value = 3078.1052
$ws.Range("J136").Value = 11692.346
$ws.Range("K136").Value = 9234.3156
$ws.Range("L136").Value = 35077.038
$ws.Range("M136").Value = -6684.3156
$ws.Range("N136").Value = -40177.038
